# Add a new "description" column header to the "params" sheet.
#
# The new column M is inserted right after the existing "ui variable"
# column (L) by copying L's formatting, so the new header cell picks up
# the same style as the rest of the header row. The copied value is then
# overwritten with the real header text, "description", which is appended
# as a new shared string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("params")

$ws.Columns("L").Copy()
$ws.Columns("M").Insert(-4161)
$ws.Range("M1").Value = "description"

# Update the active view to reflect the newly added column: scroll so
# column E becomes the left-most visible column and select the new
# header cell, matching how Excel leaves the sheet after such an edit.
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("M1").Select()
